$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (record 112438831 originally) and Row 3 (record 112438832 originally)
# swap their identifying/location/taxon data, while B column gets new distinct values.

$ws.Range("A2").Value = 112438832
$ws.Range("B2").Value = 90166
$ws.Range("E2").Value = 1339
$ws.Range("F2").Value = "Brandticka"
$ws.Range("G2").Value = "Pycnoporellus fulgens"
$ws.Range("H2").Value = "(Fr.) Donk"
$ws.Range("Q2").Value = 503830
$ws.Range("R2").Value = 6543680

$ws.Range("A3").Value = 112438831
$ws.Range("B3").Value = 89072
$ws.Range("E3").Value = 256703
$ws.Range("F3").Value = "Tallfingersvamp"
$ws.Range("G3").Value = "Ramaria eosanguinea"
$ws.Range("H3").Value = "R.H.Petersen"
$ws.Range("Q3").Value = 503890
$ws.Range("R3").Value = 6543669
